$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "[공지]2021학년도 데이터 사이언스 대학원 공식 출범"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/notice-class-of-2021-2022/#utm_source=rss&utm_medium=rss&utm_campaign=notice-class-of-2021-2022"

$ws.Range("D12").Value = "“머신러닝 파워드 애플리케이션”이 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/08/30/%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-%ed%8c%8c%ec%9b%8c%eb%93%9c-%ec%95%a0%ed%94%8c%eb%a6%ac%ec%bc%80%ec%9d%b4%ec%85%98%ec%9d%b4-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

$ws.Range("D52").Value = "파이썬 관련 아르바이트 모집"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2337&utm_source=rss&utm_medium=rss&utm_campaign=%25ed%258c%258c%25ec%259d%25b4%25ec%258d%25ac-%25ea%25b4%2580%25eb%25a0%25a8-%25ec%2595%2584%25eb%25a5%25b4%25eb%25b0%2594%25ec%259d%25b4%25ed%258a%25b8-%25eb%25aa%25a8%25ec%25a7%2591"
